$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Add "type=custom-type" line to the namespace shared-string cells ---
$ws1Cur = $ws1.Range("A1").Value2
$ws1.Range("A1").Value2 = $ws1Cur + "`n" + "type=custom-type"

$ws2Cur = $ws2.Range("A1").Value2
$ws2.Range("A1").Value2 = $ws2Cur + "`n" + "type=custom-type"

# --- Enable wrap text for the header cells (style shared by A1:C1 on sheet1, A1 on sheet2) ---
$ws1.Range("A1:C1").WrapText = $true
$ws2.Range("A1").WrapText = $true

# --- Row 1 height grows to fit the wrapped two-line text ---
$ws1.Rows.Item(1).RowHeight = 23.85
$ws2.Rows.Item(1).RowHeight = 23.85

# --- Update the saved selection on each sheet ---
$ws1.Activate()
$ws1.Range("A1").Select()

$ws2.Activate()
$ws2.Range("D16").Select()
